$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1104.375
$ws.Range("J19").Value = 1252.1666
$ws.Range("L19").Value = 1252.1666
$ws.Range("N19").Value = -1602.1666
$ws.Range("H39").Value = 237.4
$ws.Range("I39").Value = 200.375
$ws.Range("K39").Value = 601.125
$ws.Range("M39").Value = -305.125
$ws.Range("H41").Value = 1427.091
$ws.Range("I41").Value = 1915.125
$ws.Range("J41").Value = 125.666664
$ws.Range("K41").Value = 1915.125
$ws.Range("L41").Value = 125.666664
$ws.Range("M41").Value = -1475.125
$ws.Range("N41").Value = -1005.666664
$ws.Range("H70").Value = 1976.5454
$ws.Range("I70").Value = 1936.25
$ws.Range("J70").Value = 1999.5714
$ws.Range("K70").Value = 5808.75
$ws.Range("L70").Value = 5998.7142
$ws.Range("M70").Value = -5538.75
$ws.Range("N70").Value = -6538.7142
$ws.Range("H73").Value = 1976.5454
$ws.Range("I73").Value = 1936.25
$ws.Range("J73").Value = 1999.5714
$ws.Range("K73").Value = 5808.75
$ws.Range("L73").Value = 5998.7142
$ws.Range("M73").Value = -4872.75
$ws.Range("N73").Value = -7870.7142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1772.2142
$ws.Range("I2").Value = 1913.7
$ws.Range("K2").Value = 1913.7
$ws.Range("M2").Value = -1800.7
$ws.Range("H32").Value = 2794.9363
$ws.Range("I32").Value = 2365.9302
$ws.Range("K32").Value = 2365.9302
$ws.Range("M32").Value = -2078.9302
$ws.Range("H45").Value = 17192
$ws.Range("I45").Value = 12889.077
$ws.Range("K45").Value = 12889.077
$ws.Range("M45").Value = -12512.077
$ws.Range("H61").Value = 4194.6523
$ws.Range("I61").Value = 3665.389
$ws.Range("K61").Value = 3665.389
$ws.Range("M61").Value = -3453.389
$ws.Range("H116").Value = 1772.2142
$ws.Range("I116").Value = 1913.7
$ws.Range("K116").Value = 1913.7
$ws.Range("M116").Value = 380.3
$ws.Range("H132").Value = 2543.15
$ws.Range("I132").Value = 2247.75
$ws.Range("J132").Value = 3724.75
$ws.Range("K132").Value = 6743.25
$ws.Range("L132").Value = 11174.25
$ws.Range("M132").Value = -4213.25
$ws.Range("N132").Value = -16234.25
$ws.Range("H136").Value = 4194.6523
$ws.Range("I136").Value = 3665.389
$ws.Range("K136").Value = 10996.167
$ws.Range("M136").Value = -8446.167000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1772.2142
$ws.Range("I3").Value = 1913.7
$ws.Range("K3").Value = 1913.7
$ws.Range("M3").Value = -1799.7
$ws.Range("H22").Value = 575
$ws.Range("I22").Value = 435.57144
$ws.Range("K22").Value = 435.57144
$ws.Range("M22").Value = -262.57144
$ws.Range("H106").Value = 58294.668
$ws.Range("J106").Value = 58294.668
$ws.Range("L106").Value = 58294.668
$ws.Range("N106").Value = -60818.668
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H134").Value = 4612.5
$ws.Range("I134").Value = 4450
$ws.Range("J134").Value = 4666.6665
$ws.Range("K134").Value = 13350
$ws.Range("L134").Value = 13999.9995
$ws.Range("M134").Value = -10815
$ws.Range("N134").Value = -19069.9995
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 36331.668
$ws.Range("J28").Value = 36331.668
$ws.Range("L28").Value = 36331.668
$ws.Range("N28").Value = -36821.668
$ws.Range("H52").Value = 75390
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 140780
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 140780
$ws.Range("M52").Value = -9706
$ws.Range("N52").Value = -141368
$ws.Range("H94").Value = 1335.8
$ws.Range("I94").Value = 1600
$ws.Range("J94").Value = 1159.6666
$ws.Range("K94").Value = 1600
$ws.Range("L94").Value = 1159.6666
$ws.Range("M94").Value = -1149
$ws.Range("N94").Value = -2061.6666
$ws.Range("H134").Value = 25234.8
$ws.Range("I134").Value = 18916.5
$ws.Range("J134").Value = 29447
$ws.Range("K134").Value = 56749.5
$ws.Range("L134").Value = 88341
$ws.Range("M134").Value = -54214.5
$ws.Range("N134").Value = -93411
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 66
$ws.Range("I47").Value = 85
$ws.Range("J47").Value = 37.5
$ws.Range("K47").Value = 255
$ws.Range("L47").Value = 112.5
$ws.Range("M47").Value = 176
$ws.Range("N47").Value = -974.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = ""
$ws.Range("H102").Value = 3003.889
$ws.Range("I102").Value = 3003.889
$ws.Range("K102").Value = 3003.889
$ws.Range("M102").Value = -1381.889
$ws.Range("H113").Value = 4541.7334
$ws.Range("I113").Value = 4562.357
$ws.Range("J113").Value = 4523.6875
$ws.Range("K113").Value = 4562.357
$ws.Range("L113").Value = 4523.6875
$ws.Range("M113").Value = -2392.357
$ws.Range("N113").Value = -8863.6875
$ws.Range("H126").Value = 27376.615
$ws.Range("I126").Value = 37766.777
$ws.Range("J126").Value = 3998.75
$ws.Range("K126").Value = 113300.331
$ws.Range("L126").Value = 11996.25
$ws.Range("M126").Value = -110830.331
$ws.Range("N126").Value = -16936.25
$ws.Range("H132").Value = 229935.98
$ws.Range("I132").Value = 240646.88
$ws.Range("J132").Value = 5007
$ws.Range("K132").Value = 721940.64
$ws.Range("L132").Value = 15021
$ws.Range("M132").Value = -719410.64
$ws.Range("N132").Value = -20081
$ws.Range("H135").Value = 79280
$ws.Range("J135").Value = 79280
$ws.Range("L135").Value = 79280
$ws.Range("N135").Value = -89420
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 24999.666
$ws.Range("I43").Value = 24999
$ws.Range("K43").Value = 24999
$ws.Range("M43").Value = -24806
$ws.Range("H46").Value = 3053.7693
$ws.Range("I46").Value = 1977.7778
$ws.Range("J46").Value = 5474.75
$ws.Range("K46").Value = 1977.7778
$ws.Range("L46").Value = 5474.75
$ws.Range("M46").Value = -1789.7778
$ws.Range("N46").Value = -5850.75
$ws.Range("H68").Value = 4121.839
$ws.Range("I68").Value = 2718.88
$ws.Range("K68").Value = 2718.88
$ws.Range("M68").Value = -1969.88
$ws.Range("H71").Value = 4121.839
$ws.Range("I71").Value = 2718.88
$ws.Range("K71").Value = 13594.4
$ws.Range("M71").Value = -9850.400000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3076
$ws.Range("I132").Value = 3055.4167
$ws.Range("J132").Value = 3158.3333
$ws.Range("K132").Value = 9166.250100000001
$ws.Range("L132").Value = 9474.999899999999
$ws.Range("M132").Value = -6636.250100000001
$ws.Range("N132").Value = -14534.9999
